# "updated cleaning of both dataframes"
#
# Marks the last 15 rows (27-41) of the `covid_stats_df` sheet's column
# description table as "to be removed" by writing an "X" into column B
# (matching the existing red-bold / red-bold-centered formatting already
# used for other flagged rows), then updates each sheet's window state
# (selection / zoom / active tab) to match where the author ended up
# working.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("covid_stats_df")
$ws2 = $wb.Worksheets.Item("government_regulation_df")

# --- Flag rows 27-41 on covid_stats_df as "TO BE REMOVED" -----------------
# Column A gets the same red-bold font already used on the other flagged
# rows (style index 3 in the original file); column B gets an "X" marker
# using the red-bold, centered font (style index 4).
for ($r = 27; $r -le 41; $r++) {
    $colA = $ws1.Range("A$r")
    $colA.Font.Bold = $true
    $colA.Font.Italic = $false
    $colA.Font.Color = 255

    $colB = $ws1.Range("B$r")
    $colB.Value = "X"
    $colB.Font.Bold = $true
    $colB.Font.Italic = $false
    $colB.Font.Color = 255
    $colB.HorizontalAlignment = -4108   # xlCenter
}

# --- Window / selection state ---------------------------------------------
# covid_stats_df: zoomed to 85%, scrolled down, with A38 selected.
$null = $ws1.Activate()
$excel.ActiveWindow.Zoom = 85
$null = $ws1.Range("A38").Select()

# government_regulation_df: ends up the active sheet, with C7 selected.
$null = $ws2.Activate()
$null = $ws2.Range("C7").Select()
